# Regenerate save_data: replace column G (header "K", previously "Strike#")
# values with newly calculated K counts.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> new value for column G
$newValues = @{
    2  = 0
    3  = 2
    4  = 1
    5  = 0
    6  = 2
    7  = 3
    8  = 0
    9  = 1
    10 = 0
    11 = 1
    12 = 3
    13 = 1
    14 = 2
    15 = 1
    16 = 0
    17 = 2
    18 = 2
    19 = 1
    20 = 0
    21 = 1
    22 = 0
    23 = 3
    24 = 0
    25 = 3
    26 = 4
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}
